# The original (and only) sheet was called "employee" but actually holds
# address-table-shaped data. Re-purpose it as the "address" sheet and add a
# brand new "employee" sheet after it with the real employee data, modelling
# a single "has-a" relationship (employee.address_id -> address.id).

$wb = $excel.ActiveWorkbook

# --- "address" sheet (re-use the existing sheet / position) -----------------
$addressSheet = $wb.Worksheets.Item(1)
$addressSheet.Name = "address"

$addressSheet.Range("A1").Value = "id"
$addressSheet.Range("B1").Value = "house"
$addressSheet.Range("C1").Value = "street"

# "3" is a text value in the export, not a number -> force text via the
# classic leading-apostrophe trick, then drop the resulting cell style so it
# stays a plain, unstyled shared-string cell.
$addressSheet.Range("A2").Value = "'3"
$addressSheet.Range("A2").ClearFormats()
$addressSheet.Range("B2").Value = "Forest Lodge"
$addressSheet.Range("C2").Value = "S.S. Academy Road"

# --- "employee" sheet (new sheet, placed right after "address") -------------
$employeeSheet = $wb.Worksheets.Add($null, $addressSheet)
$employeeSheet.Name = "employee"

$employeeSheet.Range("A1").Value = "id"
$employeeSheet.Range("B1").Value = "name"
$employeeSheet.Range("C1").Value = "age"
$employeeSheet.Range("D1").Value = "address_id"

# Blank "id" (export has no id assigned yet for this row) but still a real
# empty shared-string cell rather than a truly blank cell.
$employeeSheet.Range("A2").Value = "'"
$employeeSheet.Range("A2").ClearFormats()

$employeeSheet.Range("B2").Value = "Mainul"

$employeeSheet.Range("C2").Value = "'25"
$employeeSheet.Range("C2").ClearFormats()

$employeeSheet.Range("D2").Value = "'3"
$employeeSheet.Range("D2").ClearFormats()

# Keep the workbook's active tab on "address" (index 0), matching the
# original activeTab="0" state.
$addressSheet.Activate()
